$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 37; $r++) {
    # --- Column B: date serial -> plain text "dd/mm/yyyy" ---
    $bCell = $ws.Cells.Item($r, 2)
    $serial = $bCell.Value2
    $dateText = [DateTime]::FromOADate($serial).ToString("dd/MM/yyyy")

    # Force text storage (otherwise a dd/mm/yyyy-looking string gets
    # silently re-parsed back into a date by the smart-entry heuristic),
    # then drop back to the workbook's default (unstyled) cell format so
    # no leftover number-format / style index is left on the cell.
    $bCell.NumberFormat = "@"
    $bCell.Value = $dateText
    $bCell.Style = "Normal"

    # --- Column C: fraction -> percentage-scale number (value * 100) ---
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value = $cCell.Value2 * 100
}
